$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old sample/test rows entirely before writing the refreshed export.
$ws.Range("A2:H3").ClearContents()

# New dataset (rows 2-8) replaces the old sample/test rows (2-3), reflecting a
# refreshed export from the DB plus newly-submitted records.
# Columns: A id, B nome, C cognome, D email, E file, F approvato, G numero_tessera, H inviato

$ids    = @(1, 2, 4, 5, 6, 7, 3)
$nomi   = @("Luca", "Giulia", "Sara", "Francesco", "giuseppe", "aaa", "Marco")
$cognomi = @("Rossi", "Bianchi", "Neri", "Gallo", "cangemi", "aaa", "Verdi")
$email  = @("luca.rossi@example.com", "giulia.bianchi@example.com", "sara.neri@example.com", "francesco.gallo@example.com", "giuseppecangemi94@gmail.com", "aaa@aaa.it", "marco.verdi@example.com")
$file   = @("file1.pdf", "file2.pdf", "file4.pdf", "file5.pdf", "Schermata 2021-09-24 alle 13.37.38.png", "Schermata 2019-06-13 alle 19.28.40.png", "file3.pdf")

for ($i = 0; $i -lt $ids.Count; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $ids[$i]
}
for ($i = 0; $i -lt $nomi.Count; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $nomi[$i]
}
for ($i = 0; $i -lt $cognomi.Count; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $cognomi[$i]
}
for ($i = 0; $i -lt $email.Count; $i++) {
    $ws.Cells.Item(2 + $i, 4).Value = $email[$i]
}
for ($i = 0; $i -lt $file.Count; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $file[$i]
}

# Only the two complete / approved records carry approvato/numero_tessera/inviato
$ws.Cells.Item(6, 6).Value = "SI"
$ws.Cells.Item(6, 7).Value = 1201
$ws.Cells.Item(6, 8).Value = "SI"

$ws.Cells.Item(8, 6).Value = "SI"
$ws.Cells.Item(8, 7).Value = 1202
$ws.Cells.Item(8, 8).Value = "SI"
